$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the TLS version comment on the "Vsds" sheet (cell A46)
# ------------------------------------------------------------------
$wsVsds = $wb.Worksheets.Item("Vsds")
$cell = $wsVsds.Range("A46")
if ($cell.Comment -ne $null) {
    $cell.Comment.Text("Deprecated TLS version to use - always set to 1.2 starting in version 6.*")
}

# ------------------------------------------------------------------
# 2. Insert a new "Router ID" row on the "Vscs" sheet before row 18
# ------------------------------------------------------------------
$wsVscs = $wb.Worksheets.Item("Vscs")
$wsVscs.Rows.Item(18).Insert()

$wsVscs.Range("A18").Value = "Router ID"
$wsVscs.Range("B18").Value = $null
$wsVscs.Range("C18").Value = $null

# Copy style from the row below (old row 18, now row 19 "VM name") to keep formatting consistent
$wsVscs.Range("A19:C19").Copy()
$wsVscs.Range("A18:C18").PasteSpecial(-4122)  # xlPasteFormats

$wsVscs.Range("A18").AddComment("The router ID of this VSC in IPV4 address format. Required when system_ip is IPV6. [default: (System IP)]")
